# Fruta / hortaliza, semanal
# Insert a new weekly data row at row 13 (pushing existing rows 13-50 down to 14-51),
# then populate the new row 13 with the week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13; this shifts rows 13..50 down to 14..51
$ws.Rows.Item(13).Insert()

# Populate the new row 13 (copy of the "template" row with updated date/price/origin)
$ws.Cells.Item(13, 1).Value = 6
$ws.Cells.Item(13, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(13, 3).Value = "Metropolitana"
$ws.Cells.Item(13, 4).Value = 44547
$ws.Cells.Item(13, 5).Value = 13
$ws.Cells.Item(13, 6).Value = "Fruta"
$ws.Cells.Item(13, 7).Value = 100101
$ws.Cells.Item(13, 8).Value = "Berries"
$ws.Cells.Item(13, 9).Value = 100101008
$ws.Cells.Item(13, 10).Value = "Mora"
$ws.Cells.Item(13, 11).Value = "Sin especificar"
$ws.Cells.Item(13, 12).Value = "Primera"
$ws.Cells.Item(13, 13).Value = 200
$ws.Cells.Item(13, 14).Value = 5000
$ws.Cells.Item(13, 15).Value = 5000
$ws.Cells.Item(13, 16).Value = 5000
$ws.Cells.Item(13, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(13, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(13, 19).Value = 2500
$ws.Cells.Item(13, 20).Value = 2
